$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BitácoraExperiencia1")

# Add the two new activity descriptions
$ws.Range("C9").Value = "Desarrollo Pagina Quienes Somos"
$ws.Range("C10").Value = "Desarrollo Pagina con Galeria de Fotos"

# Update current selection to C4 (as seen in the diff's sheetView/selection)
$ws.Range("C4").Select()

# Update the window size recorded in the workbook view (best-effort; mirrors
# the user resizing/restoring the Excel window before saving)
$win = $excel.ActiveWindow
$win.Width = 13080
$win.Height = 7905
